$wb = $excel.ActiveWorkbook

# "Norway" is the template sheet we clone for the two new markets.
$norway = $wb.Worksheets.Item("Norway")

# Create "Spain" first (so it claims the lower sheetId / later shared-string
# slots), inserted immediately after "Norway".
[void]$norway.Copy($null, $norway)
$spain = $wb.Worksheets.Item(6)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3442/T2128"
[void]$spain.Cells.Select()

# Create "Italy" second, also inserted right after "Norway" -- this puts it
# before the already-created "Spain" tab even though Italy's sheetId is
# higher (it was allocated after Spain's).
[void]$norway.Copy($null, $norway)
$italy = $wb.Worksheets.Item(6)
$italy.Name = "Italy"
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3443/T1967"
[void]$italy.Range("A11").Select()

# "Italy" becomes the active tab.
[void]$italy.Activate()
